$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Cell values: language column A switches from French (fra) to English
#    (eng), document codes are renumbered DOC001-DOC010 and the name/descr
#    columns get new English document-type labels (Sierra Leone master data).
# ---------------------------------------------------------------------------
$ws.Range("A2:A11").Value = "eng"

$ws.Range("B2").Value = "DOC001"
$ws.Range("C2").Value = "Utility bill"
$ws.Range("D2").Value = "Utility bill"

$ws.Range("B3").Value = "DOC004"
$ws.Range("C3").Value = "Attestation"
$ws.Range("D3").Value = "Attestation"

$ws.Range("B4").Value = "DOC002"
$ws.Range("C4").Value = "Birth Certificate"
$ws.Range("D4").Value = "Birth Certificate"

$ws.Range("B5").Value = "DOC005"
$ws.Range("C5").Value = "Drivers Licence"
$ws.Range("D5").Value = "Drivers Licence"

$ws.Range("B6").Value = "DOC006"
$ws.Range("C6").Value = "National ID Card"
$ws.Range("D6").Value = "National ID Card"

$ws.Range("B7").Value = "DOC007"
$ws.Range("C7").Value = "Passport"
$ws.Range("D7").Value = "Passport"

$ws.Range("B8").Value = "DOC008"
$ws.Range("C8").Value = "Voter ID"
$ws.Range("D8").Value = "Voter ID"

$ws.Range("B9").Value = "DOC009"
$ws.Range("C9").Value = "NASSIT"
$ws.Range("D9").Value = "NASSIT"

$ws.Range("B10").Value = "DOC010"
$ws.Range("C10").Value = "WASSCE"
$ws.Range("D10").Value = "WASSCE"

$ws.Range("B11").Value = "DOC003"
$ws.Range("C11").Value = "Signed Consent Form"
$ws.Range("D11").Value = "Signed Consent Form"

# ---------------------------------------------------------------------------
# 2. Fonts: Cambria / Calibri -> Arial throughout (italic state preserved per
#    column), matching the new Sierra-Leone house style.
# ---------------------------------------------------------------------------
$ws.Range("A1:E1").Font.Name = "Arial"

$ws.Range("A2:A11").Font.Name = "Arial"
$ws.Range("A2:A11").Font.Italic = $false

$ws.Range("B2:B11").Font.Name = "Arial"

$ws.Range("C2:D3").Font.Name = "Arial"
$ws.Range("C2:D3").Font.Italic = $false
$ws.Range("C5:D5").Font.Name = "Arial"
$ws.Range("C5:D5").Font.Italic = $false
$ws.Range("C7:D10").Font.Name = "Arial"
$ws.Range("C7:D10").Font.Italic = $false

$ws.Range("E2:E11").Font.Name = "Arial"

# Rows whose name/descr cells additionally pick up a white fill highlight.
$ws.Range("C4:D4").Font.Name = "Arial"
$ws.Range("C4:D4").Font.Italic = $false
$ws.Range("C4:D4").Interior.Color = 16777215

$ws.Range("C6:D6").Font.Name = "Arial"
$ws.Range("C6:D6").Font.Italic = $false
$ws.Range("C6:D6").Interior.Color = 16777215

$ws.Range("C11:D11").Font.Name = "Arial"
$ws.Range("C11:D11").Font.Italic = $false
$ws.Range("C11:D11").Interior.Color = 16777215

# ---------------------------------------------------------------------------
# 3. Row heights (content-driven re-wrap after the text/font changes).
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 25.5
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 26.5
$ws.Rows.Item(4).RowHeight = 39
$ws.Rows.Item(5).RowHeight = 26.5
$ws.Rows.Item(6).RowHeight = 26.5
$ws.Rows.Item(7).RowHeight = 15
$ws.Rows.Item(8).RowHeight = 15
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 39

# ---------------------------------------------------------------------------
# 4. Column widths - column A becomes visible/custom-sized and B:E are
#    widened to fit the new English text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.709635416666666
$ws.Columns.Item(2).ColumnWidth = 23.983072916666668
$ws.Columns.Item(3).ColumnWidth = 25.799479166666668
$ws.Columns.Item(4).ColumnWidth = 34.893229166666664
$ws.Columns.Item(5).ColumnWidth = 18.529947916666668

# ---------------------------------------------------------------------------
# 5. Selection cursor ends on B6 in the saved file.
# ---------------------------------------------------------------------------
$ws.Range("B6").Select()
